$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(23, 1).Value = 22
    $ws.Cells.Item(23, 2).NumberFormat = "@"
    $ws.Cells.Item(23, 2).Value = "2026-02-16"
    $ws.Cells.Item(23, 2).Style = "Normal"
    $ws.Cells.Item(23, 3).Value = "22:59:35"
    $ws.Cells.Item(23, 4).Value = "base_strategy"
    $ws.Cells.Item(23, 5).Value = "UP"
    $ws.Cells.Item(23, 6).Value = 0.5
    $ws.Cells.Item(23, 7).Value = ""
    $ws.Cells.Item(23, 8).Value = "OPEN"
    $ws.Cells.Item(23, 9).Value = 0
    $ws.Cells.Item(23, 10).Value = 0
    $ws.Cells.Item(23, 11).Value = 100
    $ws.Cells.Item(23, 12).Value = 0
    $ws.Cells.Item(23, 13).Value = 0
    $ws.Cells.Item(23, 14).Value = 0.6
    $ws.Cells.Item(23, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(23, 16).Value = ""
    $ws.Cells.Item(23, 17).Value = 0
}
